$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "BannedPaths" rule (row 35) is relocated further down the rule list
# and renamed "BannedPath", with its severity changed from Blocker to
# Critical and its Tags column cleared. Every row between the old and new
# position shifts up by one to fill the gap, so the cleanest way to model
# this is: delete the old row, then insert a fresh row at the new spot and
# populate it.

$ws.Rows("35").Delete()
$ws.Rows("40").Insert()

$ws.Range("A40").Value = "BannedPath"
$ws.Range("B40").Value = "Customer packages should not install content under /libs"
$ws.Range("C40").Value = "Bug"
$ws.Range("D40").Value = "Critical"

# Update the last recorded selection to reflect where the user ended up.
$ws.Range("A37").Select()
